# Corrects a batch of data-entry / computation errors in the "données06"
# sheet (columns A and C) as described in the commit message
# "modified data (because there were some logic problems)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (new A value, new C value)
$fixes = @(
    @{ Row = 18; A = 12.889999999999999;  C = 94  },
    @{ Row = 20; A = 22.009999999999998;  C = 117 },
    @{ Row = 22; A = 35.510000000000005;  C = 115 },
    @{ Row = 25; A = 9.9599999999999991;  C = 102 },
    @{ Row = 30; A = 9.85;                C = 93  },
    @{ Row = 38; A = 16.32;               C = 110 },
    @{ Row = 40; A = 22.54;               C = 116 },
    @{ Row = 45; A = 11.59;               C = 120 },
    @{ Row = 48; A = 47.65;               C = 117 },
    @{ Row = 49; A = 24.47;               C = 118 },
    @{ Row = 51; A = 56.000000000000007;  C = 122 },
    @{ Row = 52; A = 4.22;                C = 112 },
    @{ Row = 57; A = 12.64;               C = 112 },
    @{ Row = 63; A = 51.519999999999996;  C = 121 },
    @{ Row = 68; A = 11.28;               C = 106 }
)

foreach ($fix in $fixes) {
    $ws.Cells.Item($fix.Row, 1).Value = $fix.A
    $ws.Cells.Item($fix.Row, 3).Value = $fix.C
}

$wb.Save()
